$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.924.02'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '1.639.12'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").Value = "'214.72"
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = "'0.5061"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = '  +0.16%  '
$ws.Range("D8").Value = "'0.2553"
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").Value = "'0.06366"
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = "'19.47"
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("D12").Value = "'4.275"
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '1.650.66'
$ws.Range("E13").Value = '  +0.71%  '
$ws.Range("D14").Value = "'0.5439"
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").Value = '0.0₅7799'
$ws.Range("E15").Value = '  -1.76%  '
$ws.Range("D16").Value = "'64.19"
$ws.Range("E16").Value = '  -0.43%  '
$ws.Range("D17").Value = '25.951.15'
$ws.Range("E17").Value = '  -0.12%  '
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = "'196.88"
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("D20").Value = "'4.457"
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("D21").Value = "'9.929"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").Value = "'6.017"
$ws.Range("E22").Value = '  +0.37%  '
$ws.Range("D23").Value = "'1.002"
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").Value = "'1.886"
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("D25").Value = "'140.95"
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").Value = "'0.1187"
$ws.Range("E26").Value = '  +4.27%  '
$ws.Range("E27").Value = '  +0.61%  '
$ws.Range("D28").Value = "'15.69"
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("D30").Value = "'0.04937"
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").Value = "'3.252"
$ws.Range("E31").Value = '  -0.68%  '
$ws.Range("D32").Value = "'3.175"
$ws.Range("D33").Value = "'1.537"
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").Value = "'2.369"
$ws.Range("E34").Value = '  +0.33%  '
$ws.Range("D35").Value = "'0.8934"
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("E36").Value = '  -1.52%  '
$ws.Range("D37").Value = '1.132.52'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").Value = "'0.5424"
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("D39").Value = "'0.01555"
$ws.Range("E39").Value = '  -0.77%  '
$ws.Range("D40").Value = "'1.003"
$ws.Range("E40").Value = '  +0.22%  '
$ws.Range("D41").Value = "'2.543"
$ws.Range("E41").Value = '  -0.63%  '
$ws.Range("D42").Value = '0.0₈129'
$ws.Range("E42").Value = '  +9.74%  '
$ws.Range("D43").Value = "'5.577"
$ws.Range("E43").Value = '  -2.34%  '
$ws.Range("D44").Value = "'0.8155"
$ws.Range("E44").Value = '  +1.00%  '
$ws.Range("D45").Value = "'99.28"
$ws.Range("E45").Value = '  -0.46%  '
$ws.Range("D46").Value = '1.776.25'
$ws.Range("E46").Value = '  -0.26%  '
$ws.Range("D47").Value = "'0.4535"
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("D48").Value = "'1.003"
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").Value = "'54.75"
$ws.Range("E49").Value = '  -0.08%  '
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("E51").Value = '  +0.37%  '
